# Update "想去人数" (want-to-go count) and "最低票价" (lowest price) figures
# across the four sheets, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2895
$ws1.Range("F3").Value = 21297
$ws1.Range("F4").Value = 106
$ws1.Range("F5").Value = 3129
$ws1.Range("F7").Value = 620
$ws1.Range("F8").Value = 524
$ws1.Range("F9").Value = 780
$ws1.Range("F11").Value = 265
$ws1.Range("F12").Value = 74
$ws1.Range("F14").Value = 533
$ws1.Range("F16").Value = 280
$ws1.Range("F18").Value = 433
$ws1.Range("F19").Value = 85
$ws1.Range("F22").Value = 45
$ws1.Range("F23").Value = 130

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F4").Value = 126
$ws2.Range("F5").Value = 343
$ws2.Range("F14").Value = 158

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6146
$ws3.Range("F3").Value = 710
$ws3.Range("F5").Value = 1664
$ws3.Range("F6").Value = 60

# --- Sheet 4: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6146
$ws4.Range("F3").Value = 710
$ws4.Range("F5").Value = 1664
$ws4.Range("F6").Value = 2895
$ws4.Range("F7").Value = 21297
$ws4.Range("G8").Value = "不可售"
$ws4.Range("F10").Value = 106
$ws4.Range("F11").Value = 126
$ws4.Range("F12").Value = 343
$ws4.Range("F13").Value = 3129
$ws4.Range("F16").Value = 60
$ws4.Range("F17").Value = 620
$ws4.Range("F18").Value = 524
$ws4.Range("F19").Value = 780
$ws4.Range("F21").Value = 265
$ws4.Range("F23").Value = 74
$ws4.Range("F29").Value = 533
$ws4.Range("F33").Value = 280
$ws4.Range("F34").Value = 158
$ws4.Range("F35").Value = 158
$ws4.Range("F37").Value = 433
$ws4.Range("F39").Value = 85
$ws4.Range("F44").Value = 45
$ws4.Range("F50").Value = 130
